$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values for rows 2-18 (A-E, G-H); F (URL) handled via hyperlinks below ---
$ws.Range("A2").Value = "2026-02-10 19:02:44"
$ws.Range("B2").Value = "不動産管理システム開発エンジニア募集(AI駆動開発 × PoC / MVP)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("G2").Value = 405
$ws.Range("H2").Value = "🔥AI,Ai ◆開発,システム開発 ◇管理"

$ws.Range("A3").Value = "2026-02-10 19:02:44"
$ws.Range("B3").Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("G3").Value = 385
$ws.Range("H3").Value = "🔥AI,Ai ◆効率化"

$ws.Range("A4").Value = "2026-02-10 19:02:44"
$ws.Range("B4").Value = "初回 AIプロダクト開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("G4").Value = 375
$ws.Range("H4").Value = "🔥AI,Ai ◆開発"

$ws.Range("A5").Value = "2026-02-10 19:02:44"
$ws.Range("B5").Value = "【週5日】法人向け生成AIサービス(RAG・議事録機能)のコア開発を担うリードエンジニア募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("G5").Value = 375
$ws.Range("H5").Value = "🔥AI,Ai ◆開発"

$ws.Range("A6").Value = "2026-02-10 19:02:44"
$ws.Range("B6").Value = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("G6").Value = 368
$ws.Range("H6").Value = "🔥AI,Ai ◆開発"

$ws.Range("A7").Value = "2026-02-10 19:02:44"
$ws.Range("B7").Value = "企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("G7").Value = 348
$ws.Range("H7").Value = "🔥AI,Ai ◆コンサル"

$ws.Range("A8").Value = "2026-02-10 19:02:44"
$ws.Range("B8").Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("G8").Value = 310
$ws.Range("H8").Value = "🔥AI,Ai"

$ws.Range("A9").Value = "2026-02-10 19:02:44"
$ws.Range("B9").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("G9").Value = 243
$ws.Range("H9").Value = "🔥API ◆ツール"

$ws.Range("A10").Value = "2026-02-10 19:02:44"
$ws.Range("B10").Value = "【急募】新聞記事PDFをCSV・Excel化するPythonプログラム作成依頼"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("G10").Value = 198
$ws.Range("H10").Value = "🔥Python"

$ws.Range("A11").Value = "2026-02-10 19:02:44"
$ws.Range("B11").Value = "初回 自動車販売・整備業の管理システム開発"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("G11").Value = 153
$ws.Range("H11").Value = "◆開発,システム開発 ◇管理"

$ws.Range("A12").Value = "2026-02-10 19:02:44"
$ws.Range("B12").Value = "【BUYMA】商品リスト取得ツールと自動出品ツール開発のご依頼"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("G12").Value = 123
$ws.Range("H12").Value = "◆ツール,開発"

$ws.Range("A13").Value = "2026-02-10 19:02:44"
$ws.Range("B13").Value = "【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("G13").Value = 108
$ws.Range("H13").Value = "◆開発 ◇アプリ"

$ws.Range("A14").Value = "2026-02-10 19:02:44"
$ws.Range("B14").Value = "【Java/講師/福岡市内】企業向け新入社員研修のJava講師業務(サブ講師)"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("G14").Value = 78
$ws.Range("H14").Value = "★Java"

$ws.Range("A15").Value = "2026-02-10 19:02:44"
$ws.Range("B15").Value = "【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("G15").Value = 68
$ws.Range("H15").Value = "◆ツール"

$ws.Range("A16").Value = "2026-02-10 19:02:44"
$ws.Range("B16").Value = "スプレッドシート(Apps Script)で作業時間をボタン1つで計測・集計できる仕組みの開発"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("G16").Value = 68
$ws.Range("H16").Value = "◆開発"

$ws.Range("A17").Value = "2026-02-10 19:02:44"
$ws.Range("B17").Value = "【農機具管理】顧客指定で保有機情報を見れるシステム構築依頼"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = "◇管理"

$ws.Range("A18").Value = "2026-02-10 19:02:44"
$ws.Range("B18").Value = "【急募】ECサイト(WooCommerce)の決済・配送ロジックテスト、デバッグ検証依頼"
$ws.Range("C18").Value = "システム開発"
$ws.Range("D18").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E18").Value = "期限情報なし"
$ws.Range("G18").Value = 33
$ws.Range("H18").Value = "◇サイト"

# --- Column widths (raw OOXML width = ColumnWidth + 5/6 in this engine) ---
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
$ws.Columns.Item(8).ColumnWidth = 22.166666666666668

# --- Rebuild hyperlinks for column F (URL), rows 2-18 ---
$ws.Hyperlinks.Delete()
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5489563"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5489563") | Out-Null
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5423720") | Out-Null
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5489585"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5489585") | Out-Null
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5460267"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5460267") | Out-Null
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5434128"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5434128") | Out-Null
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5434363"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5434363") | Out-Null
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5427956"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5427956") | Out-Null
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5217096") | Out-Null
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5489128"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5489128") | Out-Null
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5489393"
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5489393") | Out-Null
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5489608"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5489608") | Out-Null
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5454210"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5454210") | Out-Null
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5488955"
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5488955") | Out-Null
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5489500"
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5489500") | Out-Null
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5488743"
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5488743") | Out-Null
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5489112"
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5489112") | Out-Null
$ws.Range("F18").Value = "https://www.lancers.jp/work/detail/5489409"
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5489409") | Out-Null
$ws.Range("F2:F18").Style = "Hyperlink"
